# TC30_Canine_StudyUBC-AllBreeds_StageOfDisease.xlsx
# "updated 2 icdc scripts to resolve wait time issue"
#
# The "CasesTab" row's Cypher query (cell B2 on the "startup" sheet) is
# trimmed: it no longer returns the `Cohort` column (the final
# `coalesce(co.cohort_description, '') AS `Cohort`` line is dropped, and
# the trailing comma on the now-last `Response to Treatment` line goes
# away with it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in ['T2N1M0', 'Not Applicable']  OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# The row shrinks by one wrapped line now that the Cohort column is gone
# (304.5 -> 290, matching the other two data rows).
$ws.Rows.Item(2).RowHeight = 290

# Excel leaves the cursor sitting on the cell that was just edited.
$ws.Range("B2").Select() | Out-Null
